$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 243
$ws.Range("F3").Value = 811
$ws.Range("F4").Value = 550
$ws.Range("F5").Value = 2279
$ws.Range("F6").Value = 1358
$ws.Range("F7").Value = 115
$ws.Range("F9").Value = 1142
$ws.Range("F10").Value = 33
$ws.Range("F11").Value = 2996
$ws.Range("F12").Value = 32
$ws.Range("F13").Value = 45
$ws.Range("F14").Value = 1097
$ws.Range("F15").Value = 609
$ws.Range("F16").Value = 531
$ws.Range("F17").Value = 229
$ws.Range("F18").Value = 611
$ws.Range("F19").Value = 1061
$ws.Range("F20").Value = 1061
$ws.Range("F21").Value = 141
$ws.Range("F22").Value = 522
$ws.Range("F23").Value = 169
$ws.Range("F25").Value = 207
$ws.Range("F29").Value = 831
$ws.Range("F30").Value = 65
$ws.Range("F32").Value = 1032
$ws.Range("F33").Value = 5027
$ws.Range("F34").Value = 490
$ws.Range("F35").Value = 228
$ws.Range("F36").Value = 119

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 17
$ws.Range("F21").Value = 300
$ws.Range("F22").Value = 42
$ws.Range("F27").Value = 655
$ws.Range("F30").Value = 7
$ws.Range("F36").Value = 431
$ws.Range("F37").Value = 431
$ws.Range("F43").Value = 751
$ws.Range("F44").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 419
$ws.Range("F6").Value = 400

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 243
$ws.Range("F4").Value = 419
$ws.Range("F5").Value = 811
$ws.Range("F7").Value = 550
$ws.Range("F9").Value = 2279
$ws.Range("F10").Value = 1358
$ws.Range("F11").Value = 115
$ws.Range("F14").Value = 33
$ws.Range("F16").Value = 2996
$ws.Range("F17").Value = 32
$ws.Range("F18").Value = 45
$ws.Range("F19").Value = 1097
$ws.Range("F20").Value = 609
$ws.Range("F22").Value = 400
$ws.Range("F25").Value = 531
$ws.Range("F26").Value = 229
$ws.Range("F27").Value = 1061
$ws.Range("F28").Value = 1061
$ws.Range("F29").Value = 141
$ws.Range("F31").Value = 522
$ws.Range("F32").Value = 169
$ws.Range("F33").Value = 207
$ws.Range("F34").Value = 42
$ws.Range("F38").Value = 655
$ws.Range("F39").Value = 831
$ws.Range("F40").Value = 65
$ws.Range("F41").Value = 1032
$ws.Range("F42").Value = 5027
$ws.Range("F44").Value = 490
$ws.Range("F46").Value = 431
$ws.Range("F47").Value = 228
